# Update "想去人数" (want-to-go count) figures that changed between scrapes.
# Two worksheets carry the same underlying rows: "展览" (exhibitions) and
# "全部类型" (all types, the union of all categories). Both need the same
# F-column (column 6) updates, just at different row numbers.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row -> new F value, for the "展览" sheet
$exhibitUpdates = @{
    3  = 57
    5  = 184
    6  = 9565
    7  = 859
    9  = 1209
    10 = 2058
    14 = 272
    18 = 1324
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

# Row -> new F value, for the "全部类型" sheet (rows shifted by +1 vs. 展览
# because of an extra row present only in this combined sheet)
$allUpdates = @{
    3  = 57
    6  = 184
    7  = 9565
    8  = 859
    10 = 1209
    11 = 2058
    15 = 272
    19 = 1324
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
